{"js": "const body = context.document.body;\n\n// 1) \"...to get there consequently...\" -> \"...to illuminate these mechanisms consequently...\"\nlet results = body.search(\"get there\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"illuminate these mechanisms\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Typo fix: \"underlaying\" -> \"underlying\"\nresults = body.search(\"underlaying\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"underlying\", \"Replace\");\n  await context.sync();\n}\n\n// 3) \"...that drive the behaviour that creates the overt observations.\" ->\n//    \"...that drive the behaviour that creates those overt observations.\"\nresults = body.search(\"the overt observations\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"those overt observations\", \"Replace\");\n  await context.sync();\n}\n\n// 4) Append a new sentence right after \"...those overt observations. \"\nresults = body.search(\"those overt observations. \", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"those overt observations. Different cases of Bayesian Cognitive models were implemented on the avalanche decision making data to \",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-All($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $ok = $rng.Find.Execute(\n        $findText,    # FindText\n        $false,       # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        2             # Replace (wdReplaceAll)\n    )\n    if (-not $ok) {\n        Write-Output \"NOT FOUND: $findText\"\n    }\n}\n\n# Work from the end of the paragraph backwards so each Find.Execute call resolves\n# against text that has not shifted yet.\n\n# 4) Append the new sentence right after \"...the overt observations. \" and\n#    fix \"the\" -> \"those\" at the same time.\nReplace-All \"the overt observations. \" \"those overt observations. Different cases of Bayesian Cognitive models were implemented on the avalanche decision making data to \"\n\n# 2) Typo fix: \"underlaying\" -> \"underlying\"\nReplace-All \"underlaying\" \"underlying\"\n\n# 1) \"...to get there consequently...\" -> \"...to illuminate these mechanisms consequently...\"\nReplace-All \"get there\" \"illuminate these mechanisms\"\n"}
